$wb = $excel.ActiveWorkbook

# Rename the "Clients" sheet to "Données normalisées"
$ws = $wb.Worksheets.Item("Clients")
$ws.Name = "Données normalisées"

# Make it the active/selected sheet (moves tabSelected from "Visa" to this
# sheet, and updates the workbook's active-tab pointer to the first sheet)
$ws.Activate()
